# Fix "depth" column values for rows 3-37 (A3:A37) from 0 to 0.5
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3:A37").Value = 0.5
